# Atualização de bases das ligas, do dia: 21-04-2024 às 14:32
#
# This script:
#  1) Swaps the data (columns B:AC, keeping column A / the running index)
#     between five pairs of existing rows whose match order was corrected.
#  2) Appends four new match rows (161-164), copying the row-style
#     (bold/bordered id column, date-formatted date column) from the last
#     existing data row and then filling in the actual values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($Row1, $Row2) {
    $rangeA = $ws.Range("B$($Row1):AC$($Row1)")
    $rangeB = $ws.Range("B$($Row2):AC$($Row2)")
    $valA = $rangeA.Value2
    $valB = $rangeB.Value2
    $rangeA.Value2 = $valB
    $rangeB.Value2 = $valA
}

# --- 1) Swap corrected match pairs ---
Swap-RowData 29  30
Swap-RowData 76  77
Swap-RowData 87  88
Swap-RowData 99  100
Swap-RowData 111 112

Write-Host "Swaps complete."

# --- 2) Append four new rows, matching the formatting of the last row ---
$srcFormat = $ws.Range("A160:AC160")
$dstFormat = $ws.Range("A161:AC164")
$srcFormat.Copy($dstFormat)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

$row161 = @{ A=159; B=7952752; C="Bosnia Herzegovina Premier Liga"; D="Bosnia  Herzegovina Premier Liga"; E=45402.35416666666; F="Siroki Brijeg"; G="NK Igman Konjic"; H=4; I=1; J="H"; K=1.666; L=3.75; M=4; N=1.75; O=3.5; P=3.8; Q=-0.5; R=1.775; S=2.025; T=2.5; U=1.825; V=1.975; W=0.75; X=-1; Y=-1; Z=0.7749999999999999; AA=-1; AB=0.825; AC=-1 }
$row162 = @{ A=160; B=7952459; C="Bosnia Herzegovina Premier Liga"; D="Bosnia  Herzegovina Premier Liga"; E=45402.45833333334; F="Sloga"; G="Zrinjski Mostar"; H=0; I=2; J="A"; K=3.5; L=3.5; M=1.833; N=7; O=4.2; P=1.363; Q=1.25; R=1.9; S=1.9; T=2.5; U=1.9; V=1.9; W=-1; X=-1; Y=0.363; Z=-1; AA=0.8999999999999999; AB=-1; AC=0.8999999999999999 }
$row163 = @{ A=161; B=7866183; C="Bosnia Herzegovina Premier Liga"; D="Bosnia  Herzegovina Premier Liga"; E=45402.47916666666; F="FK Drina Zvornik"; G="FK Rudar Prijedor"; H=3; I=0; J="H"; K=2.75; L=3.3; M=2.25; N=2.25; O=3.3; P=2.75; Q=-0.25; R=2.025; S=1.775; T=2.25; U=1.875; V=1.925; W=1.25; X=-1; Y=-1; Z=1.025; AA=-1; AB=0.875; AC=-1 }
$row164 = @{ A=162; B=7952753; C="Bosnia Herzegovina Premier Liga"; D="Bosnia  Herzegovina Premier Liga"; E=45402.65625; F="FK Sarajevo"; G="Borac Banja Luka"; H=1; I=1; J="D"; K=2.5; L=3; M=2.625; N=2.375; O=2.75; P=3; Q=-0.25; R=2.05; S=1.75; T=2.25; U=1.975; V=1.825; W=-1; X=1.75; Y=-1; Z=-0.5; AA=0.375; AB=-0.5; AC=0.4125 }

$newRows = @($row161, $row162, $row163, $row164)
$rowNums = @(161, 162, 163, 164)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowData = $newRows[$i]
    $r = $rowNums[$i]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = $rowData[$col]
    }
}

Write-Host "New rows complete."
